$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-149 down to 35-150
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new data record
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44742
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 100112026
$ws.Range("G34").Value = "Haba"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 21000
$ws.Range("L34").Value = 22000
$ws.Range("M34").Value = 21438
$ws.Range("N34").Value = "`$/malla 25 kilos"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 858
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
